$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 351.33334
$ws.Range("I107").Value = 579.8
$ws.Range("J107").Value = 188.14285
$ws.Range("K107").Value = 579.8
$ws.Range("L107").Value = 188.14285
$ws.Range("M107").Value = 1340.2
$ws.Range("N107").Value = -4028.14285
$ws.Range("H116").Value = 4837.2764
$ws.Range("I116").Value = 5032.7144
$ws.Range("J116").Value = 4267.25
$ws.Range("K116").Value = 5032.7144
$ws.Range("L116").Value = 4267.25
$ws.Range("M116").Value = -1590.7144
$ws.Range("N116").Value = -11151.25
$ws.Range("H137").Value = 47072.668
$ws.Range("I137").Value = 67685.336
$ws.Range("J137").Value = 12718.223
$ws.Range("K137").Value = 203056.008
$ws.Range("L137").Value = 38154.669
$ws.Range("M137").Value = -200506.008
$ws.Range("N137").Value = -43254.669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1889.1428
$ws.Range("I61").Value = 885.75
$ws.Range("K61").Value = 885.75
$ws.Range("M61").Value = -673.75
$ws.Range("H74").Value = 42822.812
$ws.Range("I74").Value = 78161.266
$ws.Range("J74").Value = 1059.1818
$ws.Range("K74").Value = 78161.266
$ws.Range("L74").Value = 1059.1818
$ws.Range("M74").Value = -77287.266
$ws.Range("N74").Value = -2807.1818
$ws.Range("H77").Value = 42822.812
$ws.Range("I77").Value = 78161.266
$ws.Range("J77").Value = 1059.1818
$ws.Range("K77").Value = 390806.33
$ws.Range("L77").Value = 5295.909000000001
$ws.Range("M77").Value = -386438.33
$ws.Range("N77").Value = -14031.909
$ws.Range("H110").Value = 684
$ws.Range("I110").Value = 528.9286
$ws.Range("K110").Value = 528.9286
$ws.Range("M110").Value = 1516.0714
$ws.Range("H122").Value = 1339.8
$ws.Range("I122").Value = 1166.3334
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 3499.0002
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -1049.0002
$ws.Range("N122").Value = -9700
$ws.Range("H132").Value = 3115923.5
$ws.Range("I132").Value = 5370117.5
$ws.Range("J132").Value = 596530.1
$ws.Range("K132").Value = 16110352.5
$ws.Range("L132").Value = 1789590.3
$ws.Range("M132").Value = -16107822.5
$ws.Range("N132").Value = -1794650.3
$ws.Range("H136").Value = 1889.1428
$ws.Range("I136").Value = 885.75
$ws.Range("K136").Value = 2657.25
$ws.Range("M136").Value = -107.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 390694.16
$ws.Range("I86").Value = 1872.2222
$ws.Range("J86").Value = 779516.1
$ws.Range("K86").Value = 1872.2222
$ws.Range("L86").Value = 779516.1
$ws.Range("M86").Value = -749.2221999999999
$ws.Range("N86").Value = -781762.1
$ws.Range("H89").Value = 390694.16
$ws.Range("I89").Value = 1872.2222
$ws.Range("J89").Value = 779516.1
$ws.Range("K89").Value = 9361.110999999999
$ws.Range("L89").Value = 3897580.5
$ws.Range("M89").Value = -3745.110999999999
$ws.Range("N89").Value = -3908812.5
$ws.Range("H134").Value = 46922.39
$ws.Range("I134").Value = 2318.4
$ws.Range("J134").Value = 81233.16
$ws.Range("K134").Value = 6955.200000000001
$ws.Range("L134").Value = 243699.48
$ws.Range("M134").Value = -4420.200000000001
$ws.Range("N134").Value = -248769.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1529.8572
$ws.Range("I7").Value = 2622.25
$ws.Range("J7").Value = 73.333336
$ws.Range("K7").Value = 2622.25
$ws.Range("L7").Value = 73.333336
$ws.Range("M7").Value = -2509.25
$ws.Range("N7").Value = -299.333336
$ws.Range("H22").Value = 2692.75
$ws.Range("I22").Value = 2692.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2692.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2342.75
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 32878.152
$ws.Range("I31").Value = 35399.277
$ws.Range("J31").Value = 14600
$ws.Range("K31").Value = 35399.277
$ws.Range("L31").Value = 14600
$ws.Range("M31").Value = -35104.277
$ws.Range("N31").Value = -15190
$ws.Range("H34").Value = 32878.152
$ws.Range("I34").Value = 35399.277
$ws.Range("J34").Value = 14600
$ws.Range("K34").Value = 35399.277
$ws.Range("L34").Value = 14600
$ws.Range("M34").Value = -35197.277
$ws.Range("N34").Value = -15004
$ws.Range("H58").Value = 3382.5386
$ws.Range("I58").Value = 1062.7693
$ws.Range("J58").Value = 5702.3076
$ws.Range("K58").Value = 1062.7693
$ws.Range("L58").Value = 5702.3076
$ws.Range("M58").Value = -859.7692999999999
$ws.Range("N58").Value = -6108.3076
$ws.Range("H132").Value = 2065.4075
$ws.Range("I132").Value = 1669.9048
$ws.Range("J132").Value = 3449.6667
$ws.Range("K132").Value = 5009.7144
$ws.Range("L132").Value = 10349.0001
$ws.Range("M132").Value = -2479.7144
$ws.Range("N132").Value = -15409.0001
$ws.Range("H134").Value = 20002176
$ws.Range("I134").Value = 1915.1818
$ws.Range("J134").Value = 35716668
$ws.Range("K134").Value = 5745.5454
$ws.Range("L134").Value = 107150004
$ws.Range("M134").Value = -3210.5454
$ws.Range("N134").Value = -107155074
$ws.Range("H136").Value = 3382.5386
$ws.Range("I136").Value = 1062.7693
$ws.Range("J136").Value = 5702.3076
$ws.Range("K136").Value = 3188.3079
$ws.Range("L136").Value = 17106.9228
$ws.Range("M136").Value = -638.3078999999998
$ws.Range("N136").Value = -22206.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3960.276
$ws.Range("I5").Value = 387.61905
$ws.Range("J5").Value = 13338.5
$ws.Range("K5").Value = 1162.85715
$ws.Range("L5").Value = 40015.5
$ws.Range("M5").Value = -1050.85715
$ws.Range("N5").Value = -40239.5
$ws.Range("H135").Value = 3960.276
$ws.Range("I135").Value = 387.61905
$ws.Range("J135").Value = 13338.5
$ws.Range("K135").Value = 3488.57145
$ws.Range("L135").Value = 120046.5
$ws.Range("M135").Value = -953.5714500000004
$ws.Range("N135").Value = -125116.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44.5
$ws.Range("I2").Value = 54.25
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 54.25
$ws.Range("L2").Value = 25
$ws.Range("M2").Value = 58.75
$ws.Range("N2").Value = -251
$ws.Range("H132").Value = 33652
$ws.Range("I132").Value = 1137.5
$ws.Range("J132").Value = 145130.28
$ws.Range("K132").Value = 3412.5
$ws.Range("L132").Value = 435390.84
$ws.Range("M132").Value = -882.5
$ws.Range("N132").Value = -440450.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 724.75
$ws.Range("I22").Value = 816.3333
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 816.3333
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -521.3333
$ws.Range("N22").Value = -1040
$ws.Range("H27").Value = 724.75
$ws.Range("I27").Value = 816.3333
$ws.Range("J27").Value = 450
$ws.Range("K27").Value = 816.3333
$ws.Range("L27").Value = 450
$ws.Range("M27").Value = -709.3333
$ws.Range("N27").Value = -664
$ws.Range("H132").Value = 302819.38
$ws.Range("I132").Value = 82221.28
$ws.Range("J132").Value = 670482.9
$ws.Range("K132").Value = 246663.84
$ws.Range("L132").Value = 2011448.7
$ws.Range("M132").Value = -244133.84
$ws.Range("N132").Value = -2016508.7
$ws.Range("H136").Value = 324523.78
$ws.Range("I136").Value = 477719.84
$ws.Range("J136").Value = 2812
$ws.Range("K136").Value = 1433159.52
$ws.Range("L136").Value = 8436
$ws.Range("M136").Value = -1430609.52
$ws.Range("N136").Value = -13536

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5765.9614
$ws.Range("I132").Value = 1704.1578
$ws.Range("J132").Value = 16790.857
$ws.Range("K132").Value = 5112.4734
$ws.Range("L132").Value = 50372.571
$ws.Range("M132").Value = -2582.4734
$ws.Range("N132").Value = -55432.571
$ws.Range("H136").Value = 3394667.2
$ws.Range("I136").Value = 3403357.2
$ws.Range("J136").Value = 3333836.8
$ws.Range("K136").Value = 10210071.6
$ws.Range("L136").Value = 10001510.4
$ws.Range("M136").Value = -10207521.6
$ws.Range("N136").Value = -10006610.4
